# "adapted cover to call" - apply the cover-letter text edits.

$d = $word.ActiveDocument

$wdFindContinue  = 1
$wdReplaceOne    = 1
$wdReplaceAll    = 2

$lq = [char]0x201C   # left curly double quote
$rq = [char]0x201D   # right curly double quote

function ReplaceText($find, $replace) {
    $d.Content.Find.Execute(
        $find, $true, $false, $false, $false, $false, $true,
        $wdFindContinue, $false, $replace, $wdReplaceAll) | Out-Null
}

# ---------------------------------------------------------------------
# 1) Title: wrap the call name in curly quotes. Only the bold title line
#    ("PLOS call: ...") is touched - the later, unquoted, italicised
#    mention of the same call name in the body text must stay as-is, so
#    we anchor the search on the unique "PLOS call: " prefix.
#    "PLOS call: Machine Learning in Health and Biomedicine"
# -> 'PLOS call: "Machine Learning in Health and Biomedicine"'
# ---------------------------------------------------------------------
ReplaceText "PLOS call: Machine Learning in Health and Biomedicine" `
    ("PLOS call: " + $lq + "Machine Learning in Health and Biomedicine" + $rq)

# ---------------------------------------------------------------------
# 2) "...even small predictive performances typically..." paragraph.
# ---------------------------------------------------------------------
ReplaceText "even small predictive performances typically" `
    "even small cross-validated predictions typically"

ReplaceText "based on the same data. Elaborating such interplay between" `
    "based on the same data. Such synthetic evidence for the interplay between"

# The stray "_GoBack" bookmark used to sit right after "reproducible
# research findings " in this paragraph; it is relocated below (step 4).

# ---------------------------------------------------------------------
# 3) "Given the fundamental nature..." paragraph.
# ---------------------------------------------------------------------
ReplaceText `
    "Given the fundamental nature of our results and conclusions, we anticipate that the manuscript should attract wide attention in " `
    "Given that our results and conclusions have far-reaching implications for health policy, we anticipate our manuscript to attract wide attention in "

# ---------------------------------------------------------------------
# 4) "We also provide executable..." paragraph.
# ---------------------------------------------------------------------
ReplaceText `
    "We also provide executable " `
    "We also provide the full programming code for our analyses and figures, executable "

ReplaceText "iPython" "Jupyter"

ReplaceText `
    " notebook, the full programming code, and an interactive Web-App to illustrate the presented findings. " `
    " notebooks with extended findings, and an interactive Web-App to illustrate the presented findings. All used data are openly available to everybody. "

# Relocate the (hidden) "_GoBack" bookmark to its new position: right
# after the newly inserted "All used data..." sentence, before
# "We hope that our work...". Bookmarks.Add with an existing name moves
# it (Word keeps bookmark names unique), so the stray one in paragraph 2
# disappears automatically.
$goBackAnchor = $d.Content.Find.Execute("We hope that our work", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
$findRange = $d.Content
$findRange.Find.Execute("We hope that our work", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$target = $d.Range($findRange.Start, $findRange.Start)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
